# Rename the three "LSP ... FWCE" sheets to "... LSP WCE" (adapt-WCE -> freq-WCE naming fix),
# matching the author's "rename adapt wce to freq wce" commit.
$wb = $excel.ActiveWorkbook

$wsSimpleLsp = $wb.Worksheets.Item("opus_base LSP Simple FWCE ")
$wsSimpleLsp.Name = "opus_base Simple LSP WCE "

$wsAonLsp = $wb.Worksheets.Item("opus_base LSP AoN FWCE")
$wsAonLsp.Name = "opus_base AoN LSP WCE"

$wsFineLsp = $wb.Worksheets.Item("opus_base LSP Fine FWCE ")
$wsFineLsp.Name = "opus_base Fine LSP WCE "

# Reflect the end-of-session navigation state captured in the saved file: the
# user ends up with the "Simple FWCE" sheet's selection moved to F23, and the
# last sheet (now "opus_base Fine LSP WCE ") active/selected at G29.
$wsSimpleFwce = $wb.Worksheets.Item("opus_base Simple FWCE")
$wsSimpleFwce.Activate()
$wsSimpleFwce.Range("F23").Select()

$wsFineLsp.Activate()
$wsFineLsp.Range("G29").Select()
